# Generate Report for handback
# Updates the "Correspond Handback DateTime" values recorded for the
# e4cbc0be... handback entry (row 3) on the zh-cn report sheet, and for
# the 1a249adb... / e4cbc0be... handback entries (row 3) on the de-de
# report sheet, reflecting the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 13:02:13"
$wsZhCn.Range("G3").Value = "2016-01-11 13:03:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 13:02:31"
$wsDeDe.Range("G3").Value = "2016-01-11 13:05:35"
